$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (D: Price, E: Volume) keep their original
# string representation instead of being auto-converted to numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '49.723.27'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.557.38'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.89%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.17'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.74'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.528'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.555'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.18'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.29'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0814'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.25'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.951.29'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.561.93'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.859'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '49.515.10'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.22'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +8.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.70'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0944'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '284.12'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.99'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.52'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.47'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.17%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -7.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.82'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.43'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.64'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.69'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.39'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0784'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.66'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.97'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.112'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.34%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '120.18'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.96%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.01'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0310'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.28'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.017.65'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.99'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +7.49%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +6.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.04'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.32'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.22'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.86%  '
